$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 24, shifting existing rows 24-55 down to 25-56.
$ws.Rows(24).Insert()

# Populate the new row 24 with this week's data point (same dimension
# members as every other row in this sub-range, new date/price/volume).
$ws.Range("A24").Value = 10
$ws.Range("B24").Value = "Vega Modelo de Temuco"
$ws.Range("C24").Value = "La Araucanía"
$ws.Range("D24").Value = 44467
$ws.Range("D24").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E24").Value = 9
$ws.Range("F24").Value = 100112035
$ws.Range("G24").Value = "Bruselas (repollito)"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 40
$ws.Range("K24").Value = 23000
$ws.Range("L24").Value = 25000
$ws.Range("M24").Value = 24000
$ws.Range("N24").Value = "$/malla 10 kilos"
$ws.Range("O24").Value = "Provincia de Quillota"
$ws.Range("P24").Value = 2400
$ws.Range("Q24").Value = 10
$ws.Range("R24").Value = "Hortaliza"
